$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3-6 observation records (columns A,B,E,F,G,H,Q,R,AC) are
# cyclically shifted up by one: row3 <- old row4, row4 <- old row5,
# row5 <- old row6, row6 <- old row3.

$data = @{
    3 = @{ A = 104729333; B = 103346; E = 221423; F = "Myskmadra"; G = "Galium odoratum"; H = "(L.) Scop."; Q = 638749.7602009142; R = 6714172.200221093; AC = "noterad" }
    4 = @{ A = 104729212; B = 108194; E = 219711; F = "Sårläka"; G = "Sanicula europaea"; H = "L."; Q = 638724.763035205; R = 6714090.951460316; AC = "noterad, fin örtskog" }
    5 = @{ A = 104729332; B = 103346; E = 221423; F = "Myskmadra"; G = "Galium odoratum"; H = "(L.) Scop."; Q = 638722.1097858821; R = 6714082.968365866; AC = "rikligt" }
    6 = @{ A = 104765665; B = 89170;  E = 3215;   F = "Rödgul trumpetsvamp"; G = "Craterellus lutescens"; H = "(Fr.) Fr."; Q = 638582.5782925152; R = 6714427.311137903; AC = "noterad" }
}

foreach ($row in 3..6) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("AC$row").Value = $vals.AC
}
